$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "27.385.41"
$ws.Range("E2").Value = "  +1.03%  "

Set-TextValue "D3" "1.821.80"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("E5").Value = "  +0.74%  "

Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  +0.10%  "

Set-TextValue "D7" "0.4455"
$ws.Range("E7").Value = "  -0.16%  "

Set-TextValue "D8" "0.3748"
$ws.Range("E8").Value = "  +1.70%  "

Set-TextValue "D9" "0.07476"
$ws.Range("E9").Value = "  +2.58%  "

Set-TextValue "D10" "0.8854"
$ws.Range("E10").Value = "  +4.40%  "

Set-TextValue "D11" "21.01"
$ws.Range("E11").Value = "  +1.06%  "

Set-TextValue "D12" "1.823.71"
$ws.Range("E12").Value = "  -0.03%  "

Set-TextValue "D13" "6.754"
$ws.Range("E13").Value = "  +1.46%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.406"
$ws.Range("E14").Value = "  +1.66%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D15" "93.33"
$ws.Range("E15").Value = "  +1.79%  "

Set-TextValue "D16" "0.07129"
$ws.Range("E16").Value = "  +0.58%  "

Set-TextValue "D17" "1.002"
$ws.Range("E17").Value = "  +0.05%  "

Set-TextValue "D18" "0.000008776"
$ws.Range("E18").Value = "  -0.22%  "

Set-TextValue "D19" "1.001"
$ws.Range("E19").Value = "  +0.07%  "

Set-TextValue "D20" "15.14"
$ws.Range("E20").Value = "  +1.39%  "

Set-TextValue "D21" "27.386.40"
$ws.Range("E21").Value = "  +1.13%  "

Set-TextValue "D22" "5.406"
$ws.Range("E22").Value = "  +4.50%  "

Set-TextValue "D23" "10.90"
$ws.Range("E23").Value = "  -0.31%  "

Set-TextValue "D24" "2.053.16"
$ws.Range("E24").Value = "  +0.49%  "

Set-TextValue "D25" "1.959"
$ws.Range("E25").Value = "  -2.18%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "151.21"
$ws.Range("E26").Value = "  -0.38%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D27" "2.331"
$ws.Range("E27").Value = "  +5.49%  "

Set-TextValue "D28" "18.56"
$ws.Range("E28").Value = "  +1.07%  "

Set-TextValue "D29" "5.353"
$ws.Range("E29").Value = "  +2.29%  "

Set-TextValue "D30" "117.71"
$ws.Range("E30").Value = "  +0.58%  "

Set-TextValue "D31" "0.08855"
$ws.Range("E31").Value = "  +0.05%  "

Set-TextValue "D32" "0.7840"
$ws.Range("E32").Value = "  +6.52%  "

Set-TextValue "D33" "1.208"
$ws.Range("E33").Value = "  +2.15%  "

Set-TextValue "D34" "4.589"
$ws.Range("E34").Value = "  +3.43%  "

Set-TextValue "D35" "2.914"
$ws.Range("E35").Value = "  -0.15%  "

Set-TextValue "D36" "1.000"
$ws.Range("E36").Value = "  +0.11%  "

Set-TextValue "D37" "1.110"
$ws.Range("E37").Value = "  +1.03%  "

Set-TextValue "D38" "0.02000"
$ws.Range("E38").Value = "  +2.71%  "

Set-TextValue "D39" "0.05298"
$ws.Range("E39").Value = "  +1.01%  "

Set-TextValue "D40" "7.302"
$ws.Range("E40").Value = "  +1.42%  "

Set-TextValue "D41" "0.5304"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D42" "2.855"
$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D43" "0.1715"
$ws.Range("E43").Value = "  +0.86%  "

Set-TextValue "D44" "2.319"
$ws.Range("E44").Value = "  +19.69%  "

Set-TextValue "D45" "8.653"
$ws.Range("E45").Value = "  +1.32%  "

Set-TextValue "D46" "0.5070"
$ws.Range("E46").Value = "  +2.55%  "

Set-TextValue "D47" "10.63"
$ws.Range("E47").Value = "  +0.26%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D48" "105.47"
$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D49" "1.694"
$ws.Range("E49").Value = "  +2.10%  "

$ws.Range("E50").Value = "  +0.12%  "

Set-TextValue "D51" "0.06400"
$ws.Range("E51").Value = "  +1.02%  "
